$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.620.47'
$ws.Range("E2").Value = '  +1.20%  '
$ws.Range("D3").Value = '3.024.72'
$ws.Range("E3").Value = '  +2.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '379.08'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.05'
$ws.Range("E6").Value = '  +1.83%  '
$ws.Range("E7").Value = '  +1.16%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.59'
$ws.Range("E10").Value = '  +1.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0861'
$ws.Range("E12").Value = '  +1.37%  '
$ws.Range("D13").Value = '3.501.72'
$ws.Range("E13").Value = '  +2.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.74'
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("D16").Value = '3.037.42'
$ws.Range("E16").Value = '  +2.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.978'
$ws.Range("E17").Value = '  -2.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.65'
$ws.Range("E18").Value = '  -10.71%  '
$ws.Range("D19").Value = '51.618.66'
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.46'
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("E22").Value = '  +1.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.94'
$ws.Range("E23").Value = '  +0.66%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.50'
$ws.Range("E24").Value = '  +0.69%  '
$ws.Range("E25").Value = '  -2.42%  '
$ws.Range("E26").Value = '  +0.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.63'
$ws.Range("E27").Value = '  +7.44%  '
$ws.Range("E28").Value = '  +5.23%  '
$ws.Range("E29").Value = '  +2.65%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.108'
$ws.Range("E31").Value = '  +0.47%  '
$ws.Range("E32").Value = '  +2.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.11'
$ws.Range("E33").Value = '  +1.62%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '50.66'
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0453'
$ws.Range("E35").Value = '  +4.55%  '
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.32'
$ws.Range("E38").Value = '  +7.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.287'
$ws.Range("E39").Value = '  +10.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.09'
$ws.Range("E40").Value = '  +3.13%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.87'
$ws.Range("E41").Value = '  +3.17%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.59'
$ws.Range("E42").Value = '  +3.44%  '
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '126.72'
$ws.Range("E44").Value = '  +7.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.73'
$ws.Range("E45").Value = '  +7.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.67'
$ws.Range("E46").Value = '  +0.92%  '
$ws.Range("E47").Value = '  +3.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.40'
$ws.Range("E48").Value = '  +3.13%  '
$ws.Range("D49").Value = '2.032.39'
$ws.Range("E49").Value = '  +1.04%  '
$ws.Range("D50").Value = '3.324.01'
$ws.Range("E50").Value = '  +2.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0322'
$ws.Range("E51").Value = '  +1.99%  '
